$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the spicule type label in A1 from "tylostyle" to "strongyl"
$ws.Range("A1").Value = "strongyl"

# Match the author's saved selection (A2 active cell)
$ws.Range("A2").Select()
